# Update cryptos list (prices / 1h volume change / a few re-ranked coins)
# as produced by the "Updated cryptos list" GitHub Actions job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "86.198.43"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "3.224.46"
$ws.Range("E3").Value = "  -1.24%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "207.94"
$ws.Range("E5").Value = "  -5.21%  "
$ws.Range("D6").Value = "618.44"
$ws.Range("E6").Value = "  -1.54%  "
$ws.Range("D7").Value = "0.356"
$ws.Range("E7").Value = "  +10.90%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.670"
$ws.Range("E8").Value = "  +13.10%  "
$ws.Range("D9").Value = "0.999"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "3.222.11"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").Value = "0.558"
$ws.Range("E11").Value = "  -6.84%  "
$ws.Range("D12").Value = "0.178"
$ws.Range("E12").Value = "  +7.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000248"
$ws.Range("E13").Value = "  -11.42%  "
$ws.Range("D14").Value = "3.828.17"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").Value = "33.37"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").Value = "5.23"
$ws.Range("E16").Value = "  -3.00%  "
$ws.Range("D17").Value = "85.981.22"
$ws.Range("E17").Value = "  +1.84%  "
$ws.Range("D18").Value = "3.225.46"
$ws.Range("E18").Value = "  -0.98%  "
$ws.Range("D19").Value = "13.76"
$ws.Range("E19").Value = "  -4.48%  "
$ws.Range("D20").Value = "2.99"
$ws.Range("E20").Value = "  -7.12%  "
$ws.Range("D21").Value = "425.23"
$ws.Range("E21").Value = "  -5.37%  "
$ws.Range("D22").Value = "8.68"
$ws.Range("E22").Value = "  -4.90%  "
$ws.Range("D23").Value = "5.24"
$ws.Range("E23").Value = "  +1.15%  "
$ws.Range("D24").Value = "7.23"
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("D25").Value = "12.17"
$ws.Range("E25").Value = "  +1.28%  "
$ws.Range("D26").Value = "4.98"
$ws.Range("E26").Value = "  -4.66%  "
$ws.Range("D27").Value = "3.414.26"
$ws.Range("E27").Value = "  +0.58%  "
$ws.Range("D28").Value = "74.75"
$ws.Range("E28").Value = "  -4.39%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000124"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.170"
$ws.Range("E31").Value = "  +9.04%  "
$ws.Range("D32").Value = "0.998"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.60"
$ws.Range("E33").Value = "  -6.47%  "
$ws.Range("D34").Value = "531.62"
$ws.Range("E34").Value = "  -6.66%  "
$ws.Range("D35").Value = "1.37"
$ws.Range("E35").Value = "  -7.97%  "
$ws.Range("D36").Value = "1.92"
$ws.Range("E36").Value = "  -4.33%  "
$ws.Range("D37").Value = "6.77"
$ws.Range("E37").Value = "  +9.53%  "
$ws.Range("D38").Value = "0.134"
$ws.Range("E38").Value = "  -13.26%  "
$ws.Range("D39").Value = "22.05"
$ws.Range("E39").Value = "  -5.16%  "
$ws.Range("E40").Value = "  +0.22%  "
$ws.Range("E41").Value = "  +2.94%  "
$ws.Range("D42").Value = "0.383"
$ws.Range("E42").Value = "  -6.17%  "
$ws.Range("D43").Value = "1.94"
$ws.Range("E43").Value = "  -5.80%  "
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("D45").Value = "154.33"
$ws.Range("E45").Value = "  -3.28%  "
$ws.Range("D46").Value = "2.84"
$ws.Range("E46").Value = "  -7.16%  "
$ws.Range("D47").Value = "175.68"
$ws.Range("E47").Value = "  -7.03%  "
$ws.Range("D48").Value = "44.07"
$ws.Range("E48").Value = "  -2.64%  "
$ws.Range("D49").Value = "1.25"
$ws.Range("E49").Value = "  -5.36%  "
$ws.Range("D50").Value = "4.13"
$ws.Range("E50").Value = "  -2.50%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "0.121"
$ws.Range("E51").Value = "  +9.23%  "
